# Apply the "Generate Report for Handoff" update:
#  - Overview sheet: update "Latest HO Xliff Generate Date" for the
#    0c613336-... row (and the other rows that shared its timestamp) from
#    2016-08-19 14:21:41 to 2016-08-19 14:21:57
#  - de-de sheet: same timestamp update in "Latest Handoff Datetime"
#    (it shared the Overview value)
#  - zh-cn sheet: update "Latest Handoff Datetime" from
#    2016-08-19 14:21:35 to 2016-08-19 14:21:52
#  - zh-cn & de-de sheets: set "Priority" to "ht" for the rows that were
#    "Ready for handoff" (rows 7, 8, 9, 10, 13, 14)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 13, 14)

foreach ($r in $rows) {
    $overview.Range("G$r").Value2 = "2016-08-19 14:21:57"
    $dede.Range("H$r").Value2 = "2016-08-19 14:21:57"
    $zhcn.Range("H$r").Value2 = "2016-08-19 14:21:52"

    $zhcn.Range("E$r").Value2 = "ht"
    $dede.Range("E$r").Value2 = "ht"
}
